# Update "想去人数" (F column) counts on each sheet to the latest scrape.
$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 505
$ws.Range("F5").Value = 1190
$ws.Range("F7").Value = 207
$ws.Range("F9").Value = 784
$ws.Range("F10").Value = 435
$ws.Range("F12").Value = 275
$ws.Range("F14").Value = 210
$ws.Range("F15").Value = 12
$ws.Range("F17").Value = 6348
$ws.Range("F20").Value = 16
$ws.Range("F21").Value = 7313
$ws.Range("F24").Value = 3320
$ws.Range("F26").Value = 827
$ws.Range("F27").Value = 4487
$ws.Range("F29").Value = 171
$ws.Range("F30").Value = 166
$ws.Range("F31").Value = 1349
$ws.Range("F32").Value = 131
$ws.Range("F33").Value = 40
$ws.Range("F35").Value = 1058
$ws.Range("F36").Value = 1431
$ws.Range("F37").Value = 2096

# 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 58

# 本地生活 (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1182
$ws.Range("F4").Value = 63

# 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1182
$ws.Range("F5").Value = 63
$ws.Range("F7").Value = 505
$ws.Range("F8").Value = 1190
$ws.Range("F10").Value = 207
$ws.Range("F12").Value = 784
$ws.Range("F13").Value = 435
$ws.Range("F15").Value = 275
$ws.Range("F16").Value = 58
$ws.Range("F18").Value = 210
$ws.Range("F19").Value = 12
$ws.Range("F21").Value = 6348
$ws.Range("F24").Value = 16
$ws.Range("F25").Value = 7313
$ws.Range("F28").Value = 3320
$ws.Range("F30").Value = 827
$ws.Range("F31").Value = 4487
$ws.Range("F34").Value = 171
$ws.Range("F35").Value = 166
$ws.Range("F36").Value = 1349
$ws.Range("F37").Value = 131
$ws.Range("F38").Value = 40
$ws.Range("F40").Value = 1058
$ws.Range("F41").Value = 1431
$ws.Range("F43").Value = 2096
